$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 8 data rows (2002年-2009年), shifting 2010年-2020年 up.
$ws.Rows.Item(2).Resize(8).EntireRow.Delete()

# Append the new 2021年 row of data (now at row 13, right after 2020年 at row 12).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 228766.77
$ws.Range("C13").Value = 248874.05
$ws.Range("D13").Value = 6057.45
$ws.Range("E13").Value = 20107.28
$ws.Range("F13").Value = 213894.93
$ws.Range("G13").Value = 24512.74

# Match the year-label formatting (bold, bordered, centered) used by the other rows.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
